# [Kadastro App] Yeni kayit eklendi: 2986
# Appends a new record (row 52) to the master "Kayitlar" sheet and to the
# filtered "Erdemli" sheet (both list the same underlying records), keeping
# every cell stored as text (matching the existing rows in these sheets).

$wb = $excel.ActiveWorkbook

$values = @(
    "2986",
    "2025-09-10",
    "Erdemli",
    "1",
    "PAYDAŞ KURUM TALEP",
    "CEMAL TİMUROĞLU (K.Teknisyeni), ENDER NUSRET ÖNAL GÜLSOY (Kontrol Memuru), SEVİL SARAÇER (Tekniker), ÖZKAN AKBAŞ (Mühendis)"
)

$sheetNames = @("Kayitlar", "Erdemli")
$newRow = 52

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6))

    # Force the incoming values to be stored as text (the sheet keeps every
    # column, including the numeric-looking ones, as text), then drop the
    # explicit "@" number format again so the new cells stay unstyled like
    # the rest of the table.
    $rowRange.NumberFormat = "@"
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($newRow, $i + 1).Value = $values[$i]
    }
    $rowRange.ClearFormats()
}
